$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Servidor de base de datos MySQL" (bullet list item) gets the
#    MySQL server version appended: " versión 8.0.32 "
# ------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Servidor de base de datos MySQL", $true, $true)
if (-not $found1) {
    throw "Could not find 'Servidor de base de datos MySQL'"
}
# Extend the matched range's text in place so the appended text
# inherits the exact same run formatting (Segoe UI, 12pt) instead of
# picking up a blank/default formatting.
$rng1.Text = "Servidor de base de datos MySQL versión 8.0.32 "

# ------------------------------------------------------------------
# 2) "MySQL Workbench" (the bullet right below it) gets the
#    Workbench version + edition appended:
#    " versión 8.0.32 (MYSQL Community Server – GPL)"
# ------------------------------------------------------------------
$rng2 = $d.Range($rng1.End, $d.Content.End)
$found2 = $rng2.Find.Execute("Workbench", $true, $true)
if (-not $found2) {
    throw "Could not find 'Workbench' list item"
}
$rng2.Text = "Workbench versión 8.0.32 (MYSQL Community Server – GPL)"
